# Add a new QA test case row ("test_AddValidItem") to the "To Do Page"
# test-scenario table, in the first still-empty row (row 13) right after
# the existing "test_DisplayListOfItems" row (row 12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = 3
$ws.Range("B13").Value = "test_AddValidItem"
$ws.Range("C13").Value = "This is to test whether users can add a To Do Item with valid parameters"
$ws.Range("D13").Value = "`"jd's test item`""
$ws.Range("E13").Value = "Item is added to the list"

# Move the active selection the way the author left it after typing the
# new row in (one row down, still in column E).
$ws.Range("E14").Select()
